# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap "Islas Malvinas" / "Nicaragua" rows (195 / 196) ---
# Row 195 currently holds "Islas Malvinas" data, row 196 holds "Nicaragua" data.
# After the edit, row 195 becomes "Nicaragua" (with refreshed figures) and
# row 196 becomes "Islas Malvinas" (keeping the old "Islas Malvinas" figures).
$ws.Range("A195").Value = "Nicaragua"
$ws.Range("B195").Value = 13
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 7
$ws.Range("E195").Value = 3
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 3

$ws.Range("A196").Value = "Islas Malvinas"
$ws.Range("B196").Value = 13
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 11
$ws.Range("E196").Value = 2
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 975798
$ws.Range("C4").Value = 15147
$ws.Range("E4").Value = 802224

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 157120
$ws.Range("C8").Value = 607
$ws.Range("E8").Value = 39224
$ws.Range("G8").Value = 19
$ws.Range("H8").Value = 5896

# --- Brasil (row 14) ---
$ws.Range("B14").Value = 59875
$ws.Range("C14").Value = 679
$ws.Range("E14").Value = 26638
$ws.Range("G14").Value = 32
$ws.Range("H14").Value = 4077

# --- Burkina Faso (row 102) ---
$ws.Range("B102").Value = 632
$ws.Range("C102").Value = 3
$ws.Range("D102").Value = 453
$ws.Range("E102").Value = 137
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 42

# --- Liberia (row 143) ---
$ws.Range("E143").Value = 83
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 12
